$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.181445333333333
$ws.Range("H2").Value = 6.544335999999999
$ws.Range("I2").Value = 0.1058843243701343
$ws.Range("J2").Value = 0.1058843243701343
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 203.0691603333333
$ws.Range("N2").Value = 609.207481
$ws.Range("O2").Value = 0.9796789863919257
$ws.Range("P2").Value = 0.9796789863919257
$ws.Range("Q2").Value = 442.9842721530684
$ws.Range("R2").Value = 3986.858449377616
$ws.Range("S2").Value = 0.103732647573727
$ws.Range("T2").Value = 0.103732647573727
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.181445333333333
$ws.Range("H3").Value = 6.544335999999999
$ws.Range("I3").Value = 0.1058843243701343
$ws.Range("J3").Value = 0.1058843243701343
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5294913333333333
$ws.Range("N3").Value = 1.588474
$ws.Range("O3").Value = 0.002554457466076205
$ws.Range("P3").Value = 0.002554457466076205
$ws.Range("Q3").Value = 1.155056398140444
$ws.Range("R3").Value = 10.395507583264
$ws.Range("S3").Value = 0.0002704770029277241
$ws.Range("T3").Value = 0.0002704770029277241
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.181445333333333
$ws.Range("H4").Value = 6.544335999999999
$ws.Range("I4").Value = 0.1058843243701343
$ws.Range("J4").Value = 0.1058843243701343
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1824963333333333
$ws.Range("N4").Value = 0.547489
$ws.Range("O4").Value = 0.0008804282371915408
$ws.Range("P4").Value = 0.0008804282371915407
$ws.Range("Q4").Value = 0.3981057747004444
$ws.Range("R4").Value = 3.582951972304
$ws.Range("S4").Value = 0.0000932235490514146
$ws.Range("T4").Value = 0.00009322354905141459
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.181445333333333
$ws.Range("H5").Value = 6.544335999999999
$ws.Range("I5").Value = 0.1058843243701343
$ws.Range("J5").Value = 0.1058843243701343
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.500179
$ws.Range("N5").Value = 10.500537
$ws.Range("O5").Value = 0.01688612790480639
$ws.Range("P5").Value = 0.01688612790480639
$ws.Range("Q5").Value = 7.635449145381331
$ws.Range("R5").Value = 68.71904230843199
$ws.Range("S5").Value = 0.001787976244428096
$ws.Range("T5").Value = 0.001787976244428096
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.793934999999999
$ws.Range("H6").Value = 26.381805
$ws.Range("I6").Value = 0.426845381730038
$ws.Range("J6").Value = 0.426845381730038
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 203.0691603333333
$ws.Range("N6").Value = 609.207481
$ws.Range("O6").Value = 0.9796789863919257
$ws.Range("P6").Value = 0.9796789863919257
$ws.Range("Q6").Value = 1785.776996475912
$ws.Range("R6").Value = 16071.99296828321
$ws.Range("S6").Value = 0.4181714509193583
$ws.Range("T6").Value = 0.4181714509193583
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.793934999999999
$ws.Range("H7").Value = 26.381805
$ws.Range("I7").Value = 0.426845381730038
$ws.Range("J7").Value = 0.426845381730038
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.5294913333333333
$ws.Range("N7").Value = 1.588474
$ws.Range("O7").Value = 0.002554457466076205
$ws.Range("P7").Value = 0.002554457466076205
$ws.Range("Q7").Value = 4.656312368396666
$ws.Range("R7").Value = 41.90681131557
$ws.Range("S7").Value = 0.001090358372220443
$ws.Range("T7").Value = 0.001090358372220443
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.793934999999999
$ws.Range("H8").Value = 26.381805
$ws.Range("I8").Value = 0.426845381730038
$ws.Range("J8").Value = 0.426845381730038
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1824963333333333
$ws.Range("N8").Value = 0.547489
$ws.Range("O8").Value = 0.0008804282371915408
$ws.Range("P8").Value = 0.0008804282371915407
$ws.Range("Q8").Value = 1.604860893071667
$ws.Range("R8").Value = 14.443748037645
$ws.Range("S8").Value = 0.0003758067269899276
$ws.Range("T8").Value = 0.0003758067269899276
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.793934999999999
$ws.Range("H9").Value = 26.381805
$ws.Range("I9").Value = 0.426845381730038
$ws.Range("J9").Value = 0.426845381730038
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.500179
$ws.Range("N9").Value = 10.500537
$ws.Range("O9").Value = 0.01688612790480639
$ws.Range("P9").Value = 0.01688612790480639
$ws.Range("Q9").Value = 30.78034661436499
$ws.Range("R9").Value = 277.023119529285
$ws.Range("S9").Value = 0.007207765711469332
$ws.Range("T9").Value = 0.007207765711469332
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.198723333333334
$ws.Range("H10").Value = 18.59617
$ws.Range("I10").Value = 0.3008774146563012
$ws.Range("J10").Value = 0.3008774146563012
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 203.0691603333333
$ws.Range("N10").Value = 609.207481
$ws.Range("O10").Value = 0.9796789863919257
$ws.Range("P10").Value = 0.9796789863919257
$ws.Range("Q10").Value = 1258.769542438641
$ws.Range("R10").Value = 11328.92588194777
$ws.Range("S10").Value = 0.2947632806187083
$ws.Range("T10").Value = 0.2947632806187083
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.198723333333334
$ws.Range("H11").Value = 18.59617
$ws.Range("I11").Value = 0.3008774146563012
$ws.Range("J11").Value = 0.3008774146563012
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.5294913333333333
$ws.Range("N11").Value = 1.588474
$ws.Range("O11").Value = 0.002554457466076205
$ws.Range("P11").Value = 0.002554457466076205
$ws.Range("Q11").Value = 3.282170282731111
$ws.Range("R11").Value = 29.53953254458
$ws.Range("S11").Value = 0.0007685785582424949
$ws.Range("T11").Value = 0.0007685785582424949
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.198723333333334
$ws.Range("H12").Value = 18.59617
$ws.Range("I12").Value = 0.3008774146563012
$ws.Range("J12").Value = 0.3008774146563012
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1824963333333333
$ws.Range("N12").Value = 0.547489
$ws.Range("O12").Value = 0.0008804282371915408
$ws.Range("P12").Value = 0.0008804282371915407
$ws.Range("Q12").Value = 1.131244279681111
$ws.Range("R12").Value = 10.18119851713
$ws.Range("S12").Value = 0.0002649009717965955
$ws.Range("T12").Value = 0.0002649009717965955
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.198723333333334
$ws.Range("H13").Value = 18.59617
$ws.Range("I13").Value = 0.3008774146563012
$ws.Range("J13").Value = 0.3008774146563012
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.500179
$ws.Range("N13").Value = 10.500537
$ws.Range("O13").Value = 0.01688612790480639
$ws.Range("P13").Value = 0.01688612790480639
$ws.Range("Q13").Value = 21.69664123814333
$ws.Range("R13").Value = 195.26977114329
$ws.Range("S13").Value = 0.005080654507553773
$ws.Range("T13").Value = 0.005080654507553773
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.428052000000001
$ws.Range("H14").Value = 10.284156
$ws.Range("I14").Value = 0.1663928792435264
$ws.Range("J14").Value = 0.1663928792435264
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 203.0691603333333
$ws.Range("N14").Value = 609.207481
$ws.Range("O14").Value = 0.9796789863919257
$ws.Range("P14").Value = 0.9796789863919257
$ws.Range("Q14").Value = 696.1316412190041
$ws.Range("R14").Value = 6265.184770971037
$ws.Range("S14").Value = 0.1630116072801321
$ws.Range("T14").Value = 0.1630116072801321
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.428052000000001
$ws.Range("H15").Value = 10.284156
$ws.Range("I15").Value = 0.1663928792435264
$ws.Range("J15").Value = 0.1663928792435264
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.5294913333333333
$ws.Range("N15").Value = 1.588474
$ws.Range("O15").Value = 0.002554457466076205
$ws.Range("P15").Value = 0.002554457466076205
$ws.Range("Q15").Value = 1.815123824216
$ws.Range("R15").Value = 16.336114417944
$ws.Range("S15").Value = 0.0004250435326855424
$ws.Range("T15").Value = 0.0004250435326855424
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.428052000000001
$ws.Range("H16").Value = 10.284156
$ws.Range("I16").Value = 0.1663928792435264
$ws.Range("J16").Value = 0.1663928792435264
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1824963333333333
$ws.Range("N16").Value = 0.547489
$ws.Range("O16").Value = 0.0008804282371915408
$ws.Range("P16").Value = 0.0008804282371915407
$ws.Range("Q16").Value = 0.6256069204760001
$ws.Range("R16").Value = 5.630462284284
$ws.Range("S16").Value = 0.0001464969893536029
$ws.Range("T16").Value = 0.0001464969893536029
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.428052000000001
$ws.Range("H17").Value = 10.284156
$ws.Range("I17").Value = 0.1663928792435264
$ws.Range("J17").Value = 0.1663928792435264
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.500179
$ws.Range("N17").Value = 10.500537
$ws.Range("O17").Value = 0.01688612790480639
$ws.Range("P17").Value = 0.01688612790480639
$ws.Range("Q17").Value = 11.998795621308
$ws.Range("R17").Value = 107.989160591772
$ws.Range("S17").Value = 0.002809731441355192
$ws.Range("T17").Value = 0.002809731441355192
